$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180, shifting existing rows 180-196 down to 181-197
$ws.Rows("180:180").Insert()

# Populate the newly inserted row 180 with the new data record
$ws.Cells.Item(180, 1).Value = 11
$ws.Cells.Item(180, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(180, 3).Value = "Bíobío"
$ws.Cells.Item(180, 4).Value = 45223
$ws.Cells.Item(180, 5).Value = 8
$ws.Cells.Item(180, 6).Value = 100112001
$ws.Cells.Item(180, 7).Value = "Berenjena"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 80
$ws.Cells.Item(180, 11).Value = 10000
$ws.Cells.Item(180, 12).Value = 10000
$ws.Cells.Item(180, 13).Value = 10000
$ws.Cells.Item(180, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(180, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(180, 16).Value = 200
$ws.Cells.Item(180, 17).Value = 50
$ws.Cells.Item(180, 18).Value = "Hortaliza"
